$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 11 new training-log rows (280-290), dated 2025-09-05 (serial 45905).
# Formats are copied from existing rows so that the reused style indices
# (date style, data style, empty/filled "localisation douleur" style)
# match the established pattern instead of creating new style entries:
#   row 259 -> template for rows where column G ("Localisation douleur") is empty
#   row 261 -> template for rows where column G has a value

# --- Row 280 ---
$ws.Range("A261:I261").Copy() | Out-Null
$ws.Range("A280:I280").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(280,1).Value = 45905
$ws.Cells.Item(280,2).Value = "Ilyes Boughanmi"
$ws.Cells.Item(280,3).Value = 60
$ws.Cells.Item(280,4).Value = 4
$ws.Cells.Item(280,5).Value = 5
$ws.Cells.Item(280,6).Value = 6
$ws.Cells.Item(280,7).Value = "Adducteur "
$ws.Cells.Item(280,8).Value = 10

# --- Row 281 ---
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A281:I281").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(281,1).Value = 45905
$ws.Cells.Item(281,2).Value = "Amir Etien"
$ws.Cells.Item(281,3).Value = 60
$ws.Cells.Item(281,4).Value = 4
$ws.Cells.Item(281,5).Value = 3
$ws.Cells.Item(281,6).Value = 0
$ws.Cells.Item(281,8).Value = 5

# --- Row 282 ---
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A282:I282").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(282,1).Value = 45905
$ws.Cells.Item(282,2).Value = "Omar Benyounes"
$ws.Cells.Item(282,3).Value = 60
$ws.Cells.Item(282,4).Value = 5
$ws.Cells.Item(282,5).Value = 0
$ws.Cells.Item(282,6).Value = 0
$ws.Cells.Item(282,8).Value = 5

# --- Row 283 ---
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A283:I283").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(283,1).Value = 45905
$ws.Cells.Item(283,2).Value = "Naim Ighbane"
$ws.Cells.Item(283,3).Value = 60
$ws.Cells.Item(283,4).Value = 4
$ws.Cells.Item(283,5).Value = 0
$ws.Cells.Item(283,6).Value = 0
$ws.Cells.Item(283,8).Value = 6

# --- Row 284 ---
$ws.Range("A261:I261").Copy() | Out-Null
$ws.Range("A284:I284").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(284,1).Value = 45905
$ws.Cells.Item(284,2).Value = "Karahali Souaré"
$ws.Cells.Item(284,3).Value = 60
$ws.Cells.Item(284,4).Value = 3
$ws.Cells.Item(284,5).Value = 5
$ws.Cells.Item(284,6).Value = 6
$ws.Cells.Item(284,7).Value = "Cheville "
$ws.Cells.Item(284,8).Value = 7

# --- Row 285 ---
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A285:I285").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(285,1).Value = 45905
$ws.Cells.Item(285,2).Value = "Jeremie Laurent"
$ws.Cells.Item(285,3).Value = 60
$ws.Cells.Item(285,4).Value = 7
$ws.Cells.Item(285,5).Value = 6
$ws.Cells.Item(285,6).Value = 0
$ws.Cells.Item(285,8).Value = 8

# --- Row 286 ---
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A286:I286").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(286,1).Value = 45905
$ws.Cells.Item(286,2).Value = "Sofiane Belle"
$ws.Cells.Item(286,3).Value = 60
$ws.Cells.Item(286,4).Value = 3
$ws.Cells.Item(286,5).Value = 3
$ws.Cells.Item(286,6).Value = 0
$ws.Cells.Item(286,8).Value = 6

# --- Row 287 ---
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A287:I287").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(287,1).Value = 45905
$ws.Cells.Item(287,2).Value = "Yoan Zouma"
$ws.Cells.Item(287,3).Value = 60
$ws.Cells.Item(287,4).Value = 3
$ws.Cells.Item(287,5).Value = 3
$ws.Cells.Item(287,6).Value = 4
$ws.Cells.Item(287,8).Value = 4

# --- Row 288 ---
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A288:I288").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(288,1).Value = 45905
$ws.Cells.Item(288,2).Value = "Ilan Ihaddadene"
$ws.Cells.Item(288,3).Value = 60
$ws.Cells.Item(288,4).Value = 6
$ws.Cells.Item(288,5).Value = 6
$ws.Cells.Item(288,6).Value = 0
$ws.Cells.Item(288,8).Value = 7

# --- Row 289 ---
$ws.Range("A259:I259").Copy() | Out-Null
$ws.Range("A289:I289").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(289,1).Value = 45905
$ws.Cells.Item(289,2).Value = "Naim Dhib"
$ws.Cells.Item(289,3).Value = 60
$ws.Cells.Item(289,4).Value = 4
$ws.Cells.Item(289,5).Value = 5
$ws.Cells.Item(289,6).Value = 0
$ws.Cells.Item(289,8).Value = 5

# --- Row 290 ---
$ws.Range("A261:I261").Copy() | Out-Null
$ws.Range("A290:I290").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(290,1).Value = 45905
$ws.Cells.Item(290,2).Value = "Levy Ndoutoume"
$ws.Cells.Item(290,3).Value = 60
$ws.Cells.Item(290,4).Value = 6
$ws.Cells.Item(290,5).Value = 5
$ws.Cells.Item(290,6).Value = 5
$ws.Cells.Item(290,7).Value = "Cheville"
$ws.Cells.Item(290,8).Value = 6

# Charge = Volume * Intensite, continuing the same shared formula pattern
# used throughout the sheet (e.g. I259:I279).
$ws.Range("I280:I290").Formula = "=C280*D280"

# Restore the view: scroll so row 258 is at the top and select K286,
# matching where the user was working after adding the new rows.
$excel.Application.Goto($ws.Range("A258"))
$ws.Range("K286").Select() | Out-Null

Write-Host "done"
